# DOMA-3100: add formatter convert to number for some columns
#
# The ticket-analytics export template renders numeric counters
# (processing / completed / canceled / deferred / closed / new_or_reopened)
# for two ticket rows. This adds the `:formatN()` Carbone.io formatter to
# those placeholders and applies an integer ("0") number format to the
# corresponding cells so the exported values render as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("C", "D", "E", "F", "G", "H")

$row2Values = @(
    "{d.tickets[i].processing:formatN()}",
    "{d.tickets[i].completed:formatN()}",
    "{d.tickets[i].canceled:formatN()}",
    "{d.tickets[i].deferred:formatN()}",
    "{d.tickets[i].closed:formatN()}",
    "{d.tickets[i].new_or_reopened:formatN()}"
)

$row3Values = @(
    "{d.tickets[i+1].processing:formatN()}",
    "{d.tickets[i+1].completed:formatN()}",
    "{d.tickets[i+1].canceled:formatN()}",
    "{d.tickets[i+1].deferred:formatN()}",
    "{d.tickets[i+1].closed:formatN()}",
    "{d.tickets[i+1].new_or_reopened:formatN()}"
)

for ($i = 0; $i -lt $columns.Length; $i++) {
    $col = $columns[$i]

    $cell2 = $ws.Range("$col" + "2")
    $cell2.Value = $row2Values[$i]
    $cell2.NumberFormat = "0"

    $cell3 = $ws.Range("$col" + "3")
    $cell3.Value = $row3Values[$i]
    $cell3.NumberFormat = "0"
}

$ws.Range("C2:H3").Select()
